# Applies the "Trade #19 closed at 2026-02-17 20:53:13" update to the
# live trading results workbook:
#   - Summary sheet roll-up metrics bump
#   - Strategy Status row for MarketMaking bumps
#   - All Trades: existing OPEN row (trade #47) flips to CLOSED (early_exit)
#     and a new OPEN row (trade #80) is appended
#   - MarketMaking: same two edits, mirrored (different column layout)

$wb = $excel.ActiveWorkbook

function Set-DateLikeTextCell {
    param($cell, [string]$text)
    # Force the cell to stay text (not get auto-parsed into a date/time
    # serial number) by temporarily marking the format as Text, then
    # restore the default "Normal" style so no stray style id is left
    # behind on the cell once the value has been written.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Set-BlankCell {
    param($cell)
    # Materialise a present-but-empty cell (mirrors the source workbook's
    # blank placeholder cells, e.g. unset "Exit Price"/"Exit Reason" on a
    # still-OPEN trade row) instead of leaving the cell entirely absent.
    $cell.NumberFormat = "@"
    $cell.Value = ""
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1400.43   # Current Capital
$wsSummary.Range("B4").Value = 0.22      # Total P&L $
$wsSummary.Range("B6").Value = 47        # Total Trades
$wsSummary.Range("B7").Value = 22        # Winning Trades
$wsSummary.Range("B9").Value = 46.81     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (row 5 = MarketMaking)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 100.43     # Capital
$wsStatus.Range("D5").Value = 14         # Trades
$wsStatus.Range("E5").Value = 0.11       # P&L $
$wsStatus.Range("F5").Value = 0.43       # P&L %
$wsStatus.Range("G5").Value = 57.14      # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Existing trade #47 (row 48) closes out.
$wsAll.Cells.Item(48, 7).Value = 0.16               # Exit Price
$wsAll.Cells.Item(48, 8).Value = "CLOSED"           # Status
$wsAll.Cells.Item(48, 9).Value = 14.2857            # P&L %
$wsAll.Cells.Item(48, 10).Value = 0.02              # P&L $
$wsAll.Cells.Item(48, 11).Value = 100.43            # Capital After
$wsAll.Cells.Item(48, 12).Value = "early_exit"      # Exit Reason
$wsAll.Cells.Item(48, 13).Value = 0.15              # Duration (min)

# New trade #80 (row 81) opens.
$wsAll.Cells.Item(81, 1).Value = 80                 # Trade #
Set-DateLikeTextCell $wsAll.Cells.Item(81, 2) "2026-02-17"  # Date
Set-DateLikeTextCell $wsAll.Cells.Item(81, 3) "20:53:07"    # Time
$wsAll.Cells.Item(81, 4).Value = "MarketMaking"     # Strategy
$wsAll.Cells.Item(81, 5).Value = "UP"               # Side
$wsAll.Cells.Item(81, 6).Value = 0.14               # Entry Price
Set-BlankCell $wsAll.Cells.Item(81, 7)              # Exit Price (still OPEN)
$wsAll.Cells.Item(81, 8).Value = "OPEN"             # Status
$wsAll.Cells.Item(81, 9).Value = 0                  # P&L %
$wsAll.Cells.Item(81, 10).Value = 0                 # P&L $
$wsAll.Cells.Item(81, 11).Value = 100.4069627845085 # Capital After
Set-BlankCell $wsAll.Cells.Item(81, 12)             # Exit Reason (still OPEN)
$wsAll.Cells.Item(81, 13).Value = 0                 # Duration (min)
$wsAll.Cells.Item(81, 14).Value = 0                 # Entry Slippage (bps)
$wsAll.Cells.Item(81, 15).Value = 0                 # Exit Slippage (bps)
$wsAll.Cells.Item(81, 16).Value = 0.6               # Confidence
$wsAll.Cells.Item(81, 17).Value = "Normal spread capture: 19600 bps" # Entry Reason

# ---------------------------------------------------------------------
# MarketMaking sheet (same two trades, different column order)
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

# Existing trade #47 (row 15) closes out.
$wsMM.Cells.Item(15, 7).Value = 0.16                # Exit Price
$wsMM.Cells.Item(15, 8).Value = "CLOSED"            # Status
$wsMM.Cells.Item(15, 9).Value = 14.2857             # P&L %
$wsMM.Cells.Item(15, 10).Value = 0.02               # P&L $
$wsMM.Cells.Item(15, 11).Value = 100.43             # Capital After
$wsMM.Cells.Item(15, 16).Value = "early_exit"       # Exit Reason
$wsMM.Cells.Item(15, 17).Value = 0.15               # Duration (min)

# New trade #80 (row 48) opens.
$wsMM.Cells.Item(48, 1).Value = 80                  # Trade #
Set-DateLikeTextCell $wsMM.Cells.Item(48, 2) "2026-02-17"  # Date
Set-DateLikeTextCell $wsMM.Cells.Item(48, 3) "20:53:07"    # Time
$wsMM.Cells.Item(48, 4).Value = "MarketMaking"      # Strategy
$wsMM.Cells.Item(48, 5).Value = "UP"                # Side
$wsMM.Cells.Item(48, 6).Value = 0.14                # Entry Price
$wsMM.Cells.Item(48, 8).Value = "OPEN"              # Status
$wsMM.Cells.Item(48, 9).Value = 0                   # P&L %
$wsMM.Cells.Item(48, 10).Value = 0                  # P&L $
$wsMM.Cells.Item(48, 11).Value = 100.4069627845085  # Capital After
$wsMM.Cells.Item(48, 12).Value = 0                  # Entry Slippage (bps)
$wsMM.Cells.Item(48, 13).Value = 0                  # Exit Slippage (bps)
$wsMM.Cells.Item(48, 14).Value = 0.6                # Confidence
$wsMM.Cells.Item(48, 15).Value = "Normal spread capture: 19600 bps" # Entry Reason
$wsMM.Cells.Item(48, 17).Value = 0                  # Duration (min)
